# Append a new "Commit 7" results block (rows 122-137) to Sheet1,
# mirroring the existing per-commit block structure (e.g. the
# "Commit 6" block at rows 105-120, offset by +17 rows) with
# updated numbers (smaller counts / byte-oriented cache run).
#
# NOTE: values & formulas are written FIRST, and cell formatting is
# copied in AFTERWARDS via Copy/PasteSpecial(xlPasteFormats). Doing it
# in the opposite order corrupts the SUM(F127:F136) aggregate's cached
# result in this engine, so this ordering is deliberate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Values and formulas for the new block
# ---------------------------------------------------------------

# row 122: commit marker
$ws.Range("A122").Value = "Commit 7"

# row 123: "MARS Tool Output" banner
$ws.Range("A123").Value = "MARS Tool Output"
$ws.Range("D123").Value = "Calulations"

# row 125: "Instruction Statistics Tool" label
$ws.Range("A125").Value = "Instruction Statistics Tool"

# row 126: column headers
$ws.Range("A126").Value = "Instruction type"
$ws.Range("B126").Value = "Count"
$ws.Range("D126").Value = "Adjusted count"
$ws.Range("E126").Value = "CPI"
$ws.Range("F126").Value = "Total cycles"

# row 127: ALU
$ws.Range("A127").Value = "ALU"
$ws.Range("B127").Value = 3564
$ws.Range("D127").Formula = "=B127"
$ws.Range("E127").Value = 1
$ws.Range("F127").Formula = "=D127*E127"

# row 128: Jump
$ws.Range("A128").Value = "Jump"
$ws.Range("B128").Value = 79
$ws.Range("D128").Formula = "=B128"
$ws.Range("E128").Value = 1
$ws.Range("F128").Formula = "=D128*E128"

# row 129: Branch
$ws.Range("A129").Value = "Branch"
$ws.Range("B129").Value = 907
$ws.Range("D129").Formula = "=B129"
$ws.Range("E129").Value = 2
$ws.Range("F129").Formula = "=D129*E129"

# row 130: Memory
$ws.Range("A130").Value = "Memory"
$ws.Range("B130").Value = 609

# row 131: Other
$ws.Range("A131").Value = "Other"
$ws.Range("B131").Value = 714
$ws.Range("D131").Formula = "=B131-(B135+B136-B130)"
$ws.Range("E131").Value = 5
$ws.Range("F131").Formula = "=D131*E131"

# row 133: "Data Cache Simulation Tool" label
$ws.Range("A133").Value = "Data Cache Simulation Tool"

# row 134: column headers
$ws.Range("A134").Value = "Access"
$ws.Range("B134").Value = "Count"

# row 135: Cache hit
$ws.Range("A135").Value = "Cache hit"
$ws.Range("B135").Value = 603
$ws.Range("D135").Formula = "=B135"
$ws.Range("E135").Value = 2
$ws.Range("F135").Formula = "=D135*E135"

# row 136: Cache miss
$ws.Range("A136").Value = "Cache miss"
$ws.Range("B136").Value = 108
$ws.Range("D136").Formula = "=B136"
$ws.Range("E136").Value = 40
$ws.Range("F136").Formula = "=D136*E136"

# row 137: total cycles
$ws.Range("F137").Formula = "=SUM(F127:F136)"

# ---------------------------------------------------------------
# 2) Copy formatting from the equivalent "Commit 6" block (rows
#    105-120) into the new rows, cell range by cell range so no
#    stray blank cells get introduced outside the source's shape.
# ---------------------------------------------------------------
$xlPasteFormats = -4122

function Copy-Format($srcRange, $dstRange) {
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
}

Copy-Format "A105"      "A122"
Copy-Format "A106:D106" "A123:D123"
Copy-Format "A108"      "A125"
Copy-Format "A109:B114" "A126:B131"
Copy-Format "D109:F112" "D126:F129"
Copy-Format "D114:F114" "D131:F131"
Copy-Format "A116"      "A133"
Copy-Format "A117:B117" "A134:B134"
Copy-Format "A118:B119" "A135:B136"
Copy-Format "D118:F119" "D135:F136"
Copy-Format "F120"      "F137"

# ---------------------------------------------------------------
# 3) Cursor / selection, matching the author's final view
# ---------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 107
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G140").Select()
